$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "98.692.57"
$ws.Cells.Item(2, 5).Value = "  -0.70%  "
$ws.Cells.Item(3, 4).Value = "3.348.94"
$ws.Cells.Item(3, 5).Value = "  -1.04%  "
$ws.Cells.Item(4, 5).Value = "  +0.04%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "259.53"
$ws.Cells.Item(5, 5).Value = "  -0.36%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "659.03"
$ws.Cells.Item(6, 5).Value = "  +4.48%  "
$ws.Cells.Item(7, 5).Value = "  +12.08%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.462"
$ws.Cells.Item(8, 5).Value = "  +17.24%  "
$ws.Cells.Item(9, 5).Value = "  +24.95%  "
$ws.Cells.Item(10, 5).Value = "  +0.04%  "
$ws.Cells.Item(11, 4).Value = "3.344.53"
$ws.Cells.Item(11, 5).Value = "  -1.04%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.210"
$ws.Cells.Item(12, 5).Value = "  +5.35%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "43.20"
$ws.Cells.Item(13, 5).Value = "  +19.77%  "
$ws.Cells.Item(14, 5).Value = "  +8.45%  "
$ws.Cells.Item(15, 4).Value = "98.477.06"
$ws.Cells.Item(15, 5).Value = "  -0.62%  "
$ws.Cells.Item(16, 4).Value = "3.983.53"
$ws.Cells.Item(16, 5).Value = "  +0.57%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "5.61"
$ws.Cells.Item(17, 5).Value = "  +1.30%  "
$ws.Cells.Item(18, 4).Value = "3.351.49"
$ws.Cells.Item(18, 5).Value = "  -0.75%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "7.55"
$ws.Cells.Item(19, 5).Value = "  +22.65%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "16.96"
$ws.Cells.Item(20, 5).Value = "  +10.25%  "
$ws.Cells.Item(21, 2).Value = "SuiNetwork"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "3.63"
$ws.Cells.Item(21, 5).Value = "  +1.66%  "
$ws.Cells.Item(22, 2).Value = "BitcoinCash"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "533.09"
$ws.Cells.Item(22, 5).Value = "  +7.89%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "10.47"
$ws.Cells.Item(23, 5).Value = "  +11.83%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "0.0000214"
$ws.Cells.Item(24, 5).Value = "  +0.93%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.434"
$ws.Cells.Item(25, 5).Value = "  +55.56%  "
$ws.Cells.Item(26, 2).Value = "NEARProtocol"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "6.28"
$ws.Cells.Item(26, 5).Value = "  +10.38%  "
$ws.Cells.Item(27, 2).Value = "Litecoin"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "101.52"
$ws.Cells.Item(27, 5).Value = "  +14.48%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "12.64"
$ws.Cells.Item(28, 5).Value = "  +5.45%  "
$ws.Cells.Item(29, 4).Value = "3.529.08"
$ws.Cells.Item(29, 5).Value = "  +0.35%  "
$ws.Cells.Item(30, 5).Value = "  +15.70%  "
$ws.Cells.Item(31, 5).Value = "  +0.11%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "11.03"
$ws.Cells.Item(32, 5).Value = "  +14.77%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.193"
$ws.Cells.Item(33, 5).Value = "  +0.11%  "
$ws.Cells.Item(34, 5).Value = "  +0.36%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "29.47"
$ws.Cells.Item(35, 5).Value = "  +5.51%  "
$ws.Cells.Item(36, 5).Value = "  +16.88%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "7.89"
$ws.Cells.Item(37, 5).Value = "  +7.02%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "2.11"
$ws.Cells.Item(38, 5).Value = "  +6.82%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.158"
$ws.Cells.Item(39, 5).Value = "  +4.77%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "529.83"
$ws.Cells.Item(40, 5).Value = "  +6.06%  "
$ws.Cells.Item(41, 5).Value = "  -0.77%  "
$ws.Cells.Item(42, 5).Value = "  +4.51%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.0430"
$ws.Cells.Item(43, 5).Value = "  +31.41%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "3.76"
$ws.Cells.Item(44, 5).Value = "  -2.07%  "
$ws.Cells.Item(45, 2).Value = "ARBITRUM"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.830"
$ws.Cells.Item(45, 5).Value = "  +5.57%  "
$ws.Cells.Item(46, 2).Value = "dogwifhat"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "3.40"
$ws.Cells.Item(46, 5).Value = "  +3.32%  "
$ws.Cells.Item(47, 5).Value = "  -0.01%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "7.96"
$ws.Cells.Item(48, 5).Value = "  +20.91%  "
$ws.Cells.Item(49, 5).Value = "  +5.84%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "5.14"
$ws.Cells.Item(50, 5).Value = "  +10.45%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "164.64"
$ws.Cells.Item(51, 5).Value = "  +2.78%  "
